$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Legs Update Sesi 1" - updated raw leg-calibration measurements (columns
# B:G) for rows 4-9. Columns H:M hold formulas referencing these inputs and
# will recalculate automatically.

# Row 4 (L1)
$ws.Range("B4").Value = 1500
$ws.Range("C4").Value = 1430
$ws.Range("D4").Value = 1580
$ws.Range("E4").Value = 950
$ws.Range("F4").Value = 1950
$ws.Range("G4").Value = 2100

# Row 5 (L2)
$ws.Range("B5").Value = 1400
$ws.Range("C5").Value = 1550
$ws.Range("D5").Value = 1600
$ws.Range("E5").Value = 1000
$ws.Range("F5").Value = 2080
$ws.Range("G5").Value = 2150

# Row 6 (L3)
$ws.Range("B6").Value = 1400
$ws.Range("C6").Value = 1300
$ws.Range("D6").Value = 1680
$ws.Range("E6").Value = 850
$ws.Range("F6").Value = 1800
$ws.Range("G6").Value = 2150

# Row 7 (R1)
$ws.Range("B7").Value = 1600
$ws.Range("C7").Value = 1500
$ws.Range("D7").Value = 1400
$ws.Range("E7").Value = 2050
$ws.Range("F7").Value = 1000
$ws.Range("G7").Value = 850

# Row 8 (R2)
$ws.Range("B8").Value = 1450
$ws.Range("C8").Value = 1550
$ws.Range("D8").Value = 1200
$ws.Range("E8").Value = 1870
$ws.Range("F8").Value = 1040
$ws.Range("G8").Value = 650

# Row 9 (R3)
$ws.Range("B9").Value = 1300
$ws.Range("C9").Value = 1450
$ws.Range("D9").Value = 1350
$ws.Range("E9").Value = 1800
$ws.Range("F9").Value = 1000
$ws.Range("G9").Value = 850

$excel.Calculate()

# Match the saved view state: scroll position and active selection.
$window = $excel.ActiveWindow
$window.ScrollRow = 2
$window.ScrollColumn = 2
$ws.Range("E17").Select()
